# Update the division problems in the document to match the new set of
# values. Each entry is unique within the document, so a simple
# Find/Replace (wdReplaceAll) for each pair is sufficient and unambiguous.

$d = $word.ActiveDocument

$replacements = @(
    @("401÷2=", "542÷2="),
    @("128÷4=", "101÷2="),
    @("969÷7=", "627÷5="),
    @("873÷7=", "321÷9="),
    @("772÷7=", "407÷6="),
    @("167÷2=", "691÷5="),
    @("442÷9=", "912÷7="),
    @("827÷4=", "787÷7="),
    @("756÷5=", "830÷6="),
    @("942÷3=", "341÷6="),
    @("375÷4=", "558÷8="),
    @("389÷5=", "445÷3="),
    @("659÷9=", "158÷9="),
    @("605÷2=", "862÷2="),
    @("892÷4=", "501÷8="),
    @("676÷7=", "608÷3="),
    @("609÷4=", "147÷5="),
    @("961÷6=", "115÷6="),
    @("141÷7=", "652÷8="),
    @("672÷4=", "725÷4="),
    @("185÷8=", "682÷2="),
    @("557÷9=", "617÷7="),
    @("182÷5=", "826÷8="),
    @("692÷6=", "799÷6="),
    @("208÷8=", "267÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
